# Auto commit at 2026-02-14  9:28:27.23
# Append the latest day's readings (2026-02-13, serial 46066) for both
# charging stations to the bottom of the daily data log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing row formatting (date format in A, 2dp currency in
# C:E, integer in F, …) down onto the two new rows before filling values,
# so the new cells reuse the same style indexes as the rows above them.
$ws.Range("A24:F25").Copy()
$ws.Range("A26:F27").PasteSpecial(-4122)

# Row 26 - 四方坪站 (Sifangping station)
$ws.Cells.Item(26, 1).Value = 46066
$ws.Cells.Item(26, 2).Value = "四方坪站"
$ws.Cells.Item(26, 3).Value = 10639.15
$ws.Cells.Item(26, 4).Value = 9705.68
$ws.Cells.Item(26, 5).Value = 3973.66
$ws.Cells.Item(26, 6).Value = 435

# Row 27 - 高岭站 (Gaoling station)
$ws.Cells.Item(27, 1).Value = 46066
$ws.Cells.Item(27, 2).Value = "高岭站"
$ws.Cells.Item(27, 3).Value = 3269.93
$ws.Cells.Item(27, 4).Value = 3028.06
$ws.Cells.Item(27, 5).Value = 886.53
$ws.Cells.Item(27, 6).Value = 106

# Scroll / selection state as left by the edit session
$ws.Range("G28").Select()
$excel.ActiveWindow.ScrollRow = 16
